$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize columns A and B to the new (equal) width used after the data refresh.
# (ColumnWidth is quantized internally to the nearest 1/6 character by the engine,
# so 14.666666666666666 is the closest input that reproduces the target OOXML width.)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# Update the refreshed values in A1:B32
$ws.Cells.Item(1, 1).Value = -0.24002744690420741
$ws.Cells.Item(1, 2).Value = 0.23966817983475153
$ws.Cells.Item(2, 1).Value = -0.13930244796758728
$ws.Cells.Item(2, 2).Value = 0.13848884348537993
$ws.Cells.Item(3, 1).Value = -0.035552154918033096
$ws.Cells.Item(3, 2).Value = 0.035500765946057555
$ws.Cells.Item(4, 1).Value = -0.097489253776263496
$ws.Cells.Item(4, 2).Value = 0.097248037770629736
$ws.Cells.Item(5, 1).Value = -0.091248038985110469
$ws.Cells.Item(5, 2).Value = 0.090786688870180932
$ws.Cells.Item(6, 1).Value = -0.060648447275621908
$ws.Cells.Item(6, 2).Value = 0.060603359376520949
$ws.Cells.Item(7, 1).Value = -0.040603360813589617
$ws.Cells.Item(7, 2).Value = 0.040529987100335063
$ws.Cells.Item(8, 1).Value = -0.020529988544763178
$ws.Cells.Item(8, 2).Value = 0.020493761643324504
$ws.Cells.Item(9, 1).Value = -0.014493762896138129
$ws.Cells.Item(9, 2).Value = 0.014467862636347917
$ws.Cells.Item(10, 1).Value = -0.0084678638931947603
$ws.Cells.Item(10, 2).Value = 0.0084697017466979219
$ws.Cells.Item(11, 1).Value = -0.0039697029825092045
$ws.Cells.Item(11, 2).Value = 0.003967394196383367
$ws.Cells.Item(12, 1).Value = 0.002032604546259531
$ws.Cells.Item(12, 2).Value = -0.0020792253076691303
$ws.Cells.Item(13, 1).Value = -0.052358499696292249
$ws.Cells.Item(13, 2).Value = 0.052251957895824042
$ws.Cells.Item(14, 1).Value = -0.040251959244563373
$ws.Cells.Item(14, 2).Value = 0.040165489787087338
$ws.Cells.Item(15, 1).Value = -0.021048940925786574
$ws.Cells.Item(15, 2).Value = 0.02102568113991321
$ws.Cells.Item(16, 1).Value = -0.015025682413853492
$ws.Cells.Item(16, 2).Value = 0.015003395654440421
$ws.Cells.Item(17, 1).Value = -0.0090033969332186103
$ws.Cells.Item(17, 2).Value = 0.0089999986779991659
$ws.Cells.Item(18, 1).Value = -0.0903657331072516
$ws.Cells.Item(18, 2).Value = 0.090288646181694077
$ws.Cells.Item(19, 1).Value = -0.081288647409323289
$ws.Cells.Item(19, 2).Value = 0.08066518127427269
$ws.Cells.Item(20, 1).Value = -0.071665182538044547
$ws.Cells.Item(20, 2).Value = 0.071532429433343481
$ws.Cells.Item(21, 1).Value = -0.062532430704665742
$ws.Cells.Item(21, 2).Value = 0.062351734190375208
$ws.Cells.Item(22, 1).Value = -0.093936689811840779
$ws.Cells.Item(22, 2).Value = 0.093628251216456349
$ws.Cells.Item(23, 1).Value = -0.084628252472212395
$ws.Cells.Item(23, 2).Value = 0.084125185485470766
$ws.Cells.Item(24, 1).Value = -0.04212518722744818
$ws.Cells.Item(24, 2).Value = 0.041999998248549275
$ws.Cells.Item(25, 1).Value = -0.024138242831330814
$ws.Cells.Item(25, 2).Value = 0.024138308909783746
$ws.Cells.Item(26, 1).Value = -0.018138310147097769
$ws.Cells.Item(26, 2).Value = 0.018137722870413597
$ws.Cells.Item(27, 1).Value = -0.012137724108513659
$ws.Cells.Item(27, 2).Value = 0.012124973627941849
$ws.Cells.Item(28, 1).Value = -0.0061249748694640616
$ws.Cells.Item(28, 2).Value = 0.0061226325129872805
$ws.Cells.Item(29, 1).Value = 0.0058773661597779636
$ws.Cells.Item(29, 2).Value = -0.005877008075522383
$ws.Cells.Item(30, 1).Value = 0.025877006635795397
$ws.Cells.Item(30, 2).Value = -0.026163820058247467
$ws.Cells.Item(31, 1).Value = -0.036330586049281521
$ws.Cells.Item(31, 2).Value = 0.036264868753528745
$ws.Cells.Item(32, 1).Value = -0.0060004600025749255
$ws.Cells.Item(32, 2).Value = 0.0059999987619150374
